$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '79.471.26'
$ws.Range('E2').Value = '  +3.99%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.202.28'
$ws.Range('E3').Value = '  +7.16%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '207.53'
$ws.Range('E5').Value = '  +3.30%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '630.50'
$ws.Range('E6').Value = '  -0.39%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('E8').Value = '  +12.82%  '

$ws.Range('E9').Value = '  +5.67%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '3.203.80'
$ws.Range('E10').Value = '  +7.25%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.592'
$ws.Range('E11').Value = '  +37.33%  '

$ws.Range('E12').Value = '  +1.69%  '

$ws.Range('E13').Value = '  +8.38%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.794.64'
$ws.Range('E14').Value = '  +7.30%  '

$ws.Range('E15').Value = '  +21.63%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '31.86'
$ws.Range('E16').Value = '  +9.56%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '79.228.28'
$ws.Range('E17').Value = '  +3.83%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.201.72'
$ws.Range('E18').Value = '  +7.07%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '14.48'
$ws.Range('E19').Value = '  +7.75%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '9.43'
$ws.Range('E20').Value = '  +5.20%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '431.46'
$ws.Range('E21').Value = '  +15.84%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.86'
$ws.Range('E22').Value = '  +26.32%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.94'
$ws.Range('E23').Value = '  +15.13%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.81'
$ws.Range('E24').Value = '  +6.20%  '

$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.365.60'
$ws.Range('E25').Value = '  +7.40%  '

$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '4.79'
$ws.Range('E26').Value = '  +10.19%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '77.14'
$ws.Range('E27').Value = '  +6.05%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '11.07'
$ws.Range('E28').Value = '  +13.11%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.01'
$ws.Range('E29').Value = '  +0.67%  '

$ws.Range('E30').Value = '  +8.70%  '

$ws.Range('E31').Value = '  +0.03%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '8.93'
$ws.Range('E32').Value = '  +7.91%  '

$ws.Range('E33').Value = '  +6.68%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '518.65'
$ws.Range('E34').Value = '  +2.07%  '

$ws.Range('E35').Value = '  +1.30%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.138'
$ws.Range('E36').Value = '  +24.17%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '22.78'
$ws.Range('E37').Value = '  +12.10%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.02%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.402'
$ws.Range('E39').Value = '  +5.74%  '

$ws.Range('B40').Value = 'Cronos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.113'
$ws.Range('E40').Value = '  +7.38%  '

$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '164.43'
$ws.Range('E41').Value = '  +0.19%  '

$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '197.02'
$ws.Range('E42').Value = '  +5.50%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '20.02'
$ws.Range('E43').Value = '  +0.18%  '

$ws.Range('E44').Value = '  -0.25%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '5.44'
$ws.Range('E45').Value = '  +10.05%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.804'
$ws.Range('E46').Value = '  +14.45%  '

$ws.Range('E47').Value = '  +8.91%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.30'
$ws.Range('E48').Value = '  +5.58%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '43.04'
$ws.Range('E49').Value = '  +1.49%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.55'
$ws.Range('E50').Value = '  +10.27%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.628'
$ws.Range('E51').Value = '  +6.48%  '
